# 876. Middle of linked List- Java
# Record the "Java,Python" language marker against the Middle-of-Linked-List
# row (Leetcode 876), give the new "Language" header cell (C1) its own
# bold/vertical-top style, enlarge the header row so the bold label fits,
# widen column A to fit the longer GFG/numeric labels and move the active
# selection the way the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C12: "876 - Middle of the Linked List" is implemented in Java -> tag it
# (reuses the existing "Java,Python" shared string already used by C16/C17)
$ws.Range("C12").Value = "Java,Python"

# C1 ("Language" header) gets a new style: bold font + vertical-top alignment
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").VerticalAlignment = -4160

# Header row is taller to accommodate the new formatting
$ws.Rows(1).RowHeight = 36.75

# Column A is widened
$ws.Columns(1).ColumnWidth = 17

# Final selection/view state left by the author
[void]$ws.Range("C7").Select()
